# Updates cryptos list: price (D) and volume % (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to be treated as literal text so Excel does not
    # auto-convert numeric-looking strings (e.g. "543.50") into numbers
    # and lose formatting / introduce floating point drift; then restore
    # the original cell style so no formatting is left behind.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '59.336.09'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '2.518.88'
$ws.Range("E3").Value = '  +3.14%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue $ws.Range("D5") '543.50'
$ws.Range("E5").Value = '  +0.92%  '
Set-TextValue $ws.Range("D6") '144.79'
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("E8").Value = '  +0.74%  '
$ws.Range("D9").Value = '2.547.42'
$ws.Range("E9").Value = '  +3.76%  '
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("E12").Value = '  +4.63%  '
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").Value = '2.964.53'
Set-TextValue $ws.Range("D15") '23.84'
$ws.Range("D16").Value = '59.270.19'
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("E17").Value = '  +2.71%  '
$ws.Range("D18").Value = '2.536.72'
$ws.Range("E18").Value = '  +1.13%  '
Set-TextValue $ws.Range("D19") '11.28'
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("E20").Value = '  -0.84%  '
Set-TextValue $ws.Range("D21") '326.96'
$ws.Range("E21").Value = '  +1.18%  '
Set-TextValue $ws.Range("D22") '0.996'
$ws.Range("E22").Value = '  +2.90%  '
Set-TextValue $ws.Range("D23") '5.86'
$ws.Range("E23").Value = '  +2.99%  '
Set-TextValue $ws.Range("D24") '62.09'
$ws.Range("E24").Value = '  +2.45%  '
$ws.Range("E25").Value = '  -2.52%  '
$ws.Range("E26").Value = '  +2.51%  '
$ws.Range("E27").Value = '  +1.65%  '
Set-TextValue $ws.Range("D28") '8.07'
$ws.Range("E28").Value = '  +5.59%  '
Set-TextValue $ws.Range("D29") '6.92'
$ws.Range("E29").Value = '  +4.13%  '
$ws.Range("D30").Value = '0.0₃0786'
$ws.Range("E30").Value = '  +2.62%  '
Set-TextValue $ws.Range("D31") '1.84'
$ws.Range("E31").Value = '  +0.94%  '
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("E33").Value = '  +9.61%  '
$ws.Range("E34").Value = '  -0.12%  '
Set-TextValue $ws.Range("D35") '157.06'
$ws.Range("E35").Value = '  +0.54%  '
Set-TextValue $ws.Range("D36") '18.71'
$ws.Range("E36").Value = '  +1.65%  '
Set-TextValue $ws.Range("D37") '4.41'
$ws.Range("E37").Value = '  -0.85%  '
$ws.Range("E38").Value = '  -4.59%  '
Set-TextValue $ws.Range("D39") '5.65'
$ws.Range("E39").Value = '  -3.23%  '
Set-TextValue $ws.Range("D40") '36.92'
$ws.Range("E40").Value = '  +2.09%  '
Set-TextValue $ws.Range("D41") '299.77'
$ws.Range("E41").Value = '  -4.22%  '
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E44").Value = '  -0.43%  '
Set-TextValue $ws.Range("D45") '0.605'
$ws.Range("E45").Value = '  +4.36%  '
Set-TextValue $ws.Range("D46") '10.79'
$ws.Range("E46").Value = '  +0.42%  '
Set-TextValue $ws.Range("D47") '0.0937'
$ws.Range("E47").Value = '  -0.48%  '
Set-TextValue $ws.Range("D48") '18.86'
$ws.Range("E48").Value = '  +2.91%  '
Set-TextValue $ws.Range("D49") '123.99'
$ws.Range("E49").Value = '  +1.83%  '
$ws.Range("E50").Value = '  +0.11%  '
Set-TextValue $ws.Range("D51") '0.0516'
$ws.Range("E51").Value = '  -1.48%  '
